$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Англ. Строка (English) column for the two new entries
$ws.Range("C13").Value = " I went out exploring recently\nand came back with a bunch of items."
$ws.Range("C14").Value = " But at least I can recycle all\nthese extra items!"

# Название файла в скриптах и цвет column
$ws.Range("A13").Value = "SCRIPT/P01P04A/um1401.ssb"

# Переведенная строка (Russian translation) column
$ws.Range("D13").Value = " Недавно я ходил на вылазку и\nпринёс много предметов."
$ws.Range("D14").Value = " Но по крайней мере я всегда\nмогу переработать излишки!"

# Конвертированная строка column
$ws.Range("E13").Value = " Îåäàâîï ÿ öïäéì îà âúìàèëô é\nðñéîæò íîïãï ðñåäíåóïâ."
$ws.Range("E14").Value = " Îï ðï ëñàêîåê íåñå ÿ âòåãäà\níïãô ðåñåñàáïóàóû éèìéšëé!"

$ws.Range("A14").Value = "SCRIPT/P01P04A/um1501.ssb"
$ws.Range("A15").Value = "SCRIPT/P01P04A/um1601.ssb"

# Номер строки column
$ws.Range("B13").Value = 383
$ws.Range("B14").Value = 386

# These new rows wrap onto 3 lines each, same as the other multi-line entries above
$ws.Rows("13:15").RowHeight = 43.2

# Match the new active selection left behind by this edit
$ws.Range("C14").Select() | Out-Null
